$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I4").Value = 120.333336
$ws.Range("K4").Value = 120.333336
$ws.Range("M4").Value = -6.333336000000003
$ws.Range("H4").Value = 118.53846
$ws.Range("H33").Value = 566.7143
$ws.Range("L33").Value = 1083
$ws.Range("J33").Value = 1083
$ws.Range("N33").Value = -1541
$ws.Range("I38").Value = 98.5
$ws.Range("K38").Value = 295.5
$ws.Range("M38").Value = 76.5
$ws.Range("H38").Value = 98.5
$ws.Range("H52").Value = 1037.5
$ws.Range("K52").Value = 1155
$ws.Range("I52").Value = 385
$ws.Range("M52").Value = -995
$ws.Range("J58").Value = 3556.8572
$ws.Range("H58").Value = 2137.75
$ws.Range("N58").Value = -10970.5716
$ws.Range("L58").Value = 10670.5716
$ws.Range("K92").Value = 1617.1818
$ws.Range("N92").Value = -3329.3333
$ws.Range("M92").Value = -369.1818000000001
$ws.Range("I92").Value = 1617.1818
$ws.Range("L92").Value = 833.3333
$ws.Range("J92").Value = 833.3333
$ws.Range("H92").Value = 1449.2142
$ws.Range("M99").Value = -335
$ws.Range("H99").Value = 1343.5714
$ws.Range("I99").Value = 611
$ws.Range("J99").Value = 1636.6
$ws.Range("N99").Value = -7905.799999999999
$ws.Range("K99").Value = 1833
$ws.Range("L99").Value = 4909.799999999999
$ws.Range("L127").Value = 0
$ws.Range("H127").Value = 1227.5
$ws.Range("N127").ClearContents()
$ws.Range("J127").Value = 0
$ws.Range("H132").Value = 2009.1323
$ws.Range("M132").Value = -2819.4218
$ws.Range("K132").Value = 5349.4218
$ws.Range("I132").Value = 1783.1406
$ws.Range("H135").Value = 1692.2273
$ws.Range("M135").Value = -12146.6469
$ws.Range("K135").Value = 14681.6469
$ws.Range("I135").Value = 1631.2941
$ws.Range("L137").Value = 13918.845
$ws.Range("H137").Value = 3326.2307
$ws.Range("J137").Value = 4639.615
$ws.Range("N137").Value = -19018.845
$ws.Range("I137").Value = 2012.8462
$ws.Range("M137").Value = -3488.5386
$ws.Range("K137").Value = 6038.5386

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("L32").Value = 5068.5713
$ws.Range("M32").Value = -547.6389
$ws.Range("J32").Value = 5068.5713
$ws.Range("I32").Value = 834.6389
$ws.Range("N32").Value = -5642.5713
$ws.Range("K32").Value = 834.6389
$ws.Range("H32").Value = 1523.8837
$ws.Range("M45").Value = -907.4000000000001
$ws.Range("H45").Value = 1805.7307
$ws.Range("I45").Value = 1284.4
$ws.Range("K45").Value = 1284.4
$ws.Range("I61").Value = 8914.333000000001
$ws.Range("H61").Value = 11325.6
$ws.Range("M61").Value = -8702.333000000001
$ws.Range("L61").Value = 17956.584
$ws.Range("K61").Value = 8914.333000000001
$ws.Range("N61").Value = -18380.584
$ws.Range("J61").Value = 17956.584
$ws.Range("K74").Value = 1447.6957
$ws.Range("I74").Value = 1447.6957
$ws.Range("J74").Value = 2675.5
$ws.Range("M74").Value = -573.6957
$ws.Range("H74").Value = 1701.7241
$ws.Range("L74").Value = 2675.5
$ws.Range("N74").Value = -4423.5
$ws.Range("M77").Value = -2870.4785
$ws.Range("N77").Value = -22113.5
$ws.Range("H77").Value = 1701.7241
$ws.Range("I77").Value = 1447.6957
$ws.Range("L77").Value = 13377.5
$ws.Range("K77").Value = 7238.4785
$ws.Range("J77").Value = 2675.5
$ws.Range("N95").Value = -21991.5
$ws.Range("L95").Value = 16499.5
$ws.Range("H95").Value = 16499.5
$ws.Range("J95").Value = 16499.5
$ws.Range("I97").Value = 7492.8823
$ws.Range("M97").Value = -6996.8823
$ws.Range("H97").Value = 7625.9565
$ws.Range("K97").Value = 7492.8823
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("L101").Value = 0
$ws.Range("H102").Value = 1161.1052
$ws.Range("K102").Value = 1192.2778
$ws.Range("I102").Value = 1192.2778
$ws.Range("M102").Value = 429.7221999999999
$ws.Range("M110").Value = -89.5
$ws.Range("J110").Value = 2039.6
$ws.Range("N110").Value = -6129.6
$ws.Range("I110").Value = 2134.5
$ws.Range("K110").Value = 2134.5
$ws.Range("L110").Value = 2039.6
$ws.Range("H110").Value = 2119.1936
$ws.Range("J112").Value = 20057
$ws.Range("N112").Value = -23011
$ws.Range("H112").Value = 20057
$ws.Range("L112").Value = 20057
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("N125").Value = -105908.836
$ws.Range("H125").Value = 96068.836
$ws.Range("L125").Value = 96068.836
$ws.Range("J125").Value = 96068.836
$ws.Range("H132").Value = 3109.6667
$ws.Range("M132").Value = -2067.7559
$ws.Range("K132").Value = 4597.7559
$ws.Range("I132").Value = 1532.5853
$ws.Range("N136").Value = -58969.75199999999
$ws.Range("L136").Value = 53869.75199999999
$ws.Range("M136").Value = -24192.999
$ws.Range("I136").Value = 8914.333000000001
$ws.Range("H136").Value = 11325.6
$ws.Range("K136").Value = 26742.999
$ws.Range("J136").Value = 17956.584

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K11").Value = 0
$ws.Range("H11").Value = 15000
$ws.Range("M11").ClearContents()
$ws.Range("I11").Value = 0
$ws.Range("L12").Value = 750000100
$ws.Range("J12").Value = 750000100
$ws.Range("H12").Value = 176470990
$ws.Range("K12").Value = 493.84616
$ws.Range("N12").Value = -750000436
$ws.Range("M12").Value = -325.84616
$ws.Range("I12").Value = 493.84616
$ws.Range("I20").Value = 2048.5
$ws.Range("H20").Value = 2213.5833
$ws.Range("M20").Value = -1801.5
$ws.Range("K20").Value = 2048.5
$ws.Range("L43").Value = 277222
$ws.Range("N43").Value = -277584
$ws.Range("J43").Value = 277222
$ws.Range("H43").Value = 277222
$ws.Range("L86").Value = 43666.668
$ws.Range("J86").Value = 43666.668
$ws.Range("I86").Value = 2554.7144
$ws.Range("M86").Value = -1431.7144
$ws.Range("N86").Value = -45912.668
$ws.Range("H86").Value = 14888.3
$ws.Range("K86").Value = 2554.7144
$ws.Range("H89").Value = 14888.3
$ws.Range("L89").Value = 218333.34
$ws.Range("N89").Value = -229565.34
$ws.Range("I89").Value = 2554.7144
$ws.Range("M89").Value = -7157.572
$ws.Range("K89").Value = 12773.572
$ws.Range("J89").Value = 43666.668
$ws.Range("M94").Value = -20012979
$ws.Range("I94").Value = 20013430
$ws.Range("H94").Value = 11141010
$ws.Range("K94").Value = 20013430
$ws.Range("H107").Value = 2222.2307
$ws.Range("L107").Value = 1500
$ws.Range("N107").Value = -5340
$ws.Range("J107").Value = 1500

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K16").Value = 953.8461
$ws.Range("J16").Value = 499.25
$ws.Range("M16").Value = -666.8461
$ws.Range("N16").Value = -1073.25
$ws.Range("H16").Value = 846.8823
$ws.Range("L16").Value = 499.25
$ws.Range("I16").Value = 953.8461
$ws.Range("N22").Value = -2699
$ws.Range("I22").Value = 889.1667
$ws.Range("J22").Value = 1999
$ws.Range("M22").Value = -539.1667
$ws.Range("H22").Value = 1259.1111
$ws.Range("L22").Value = 1999
$ws.Range("K22").Value = 889.1667
$ws.Range("I31").Value = 1390.6428
$ws.Range("K31").Value = 1390.6428
$ws.Range("M31").Value = -1095.6428
$ws.Range("N31").Value = -3720.6296
$ws.Range("H31").Value = 2536.4878
$ws.Range("L31").Value = 3130.6296
$ws.Range("J31").Value = 3130.6296
$ws.Range("N34").Value = -3534.6296
$ws.Range("J34").Value = 3130.6296
$ws.Range("L34").Value = 3130.6296
$ws.Range("H34").Value = 2536.4878
$ws.Range("K34").Value = 1390.6428
$ws.Range("I34").Value = 1390.6428
$ws.Range("M34").Value = -1188.6428
$ws.Range("K62").Value = 3579.375
$ws.Range("M62").Value = -2955.375
$ws.Range("H62").Value = 11890.083
$ws.Range("I62").Value = 3579.375
$ws.Range("L62").Value = 28511.5
$ws.Range("J62").Value = 28511.5
$ws.Range("N62").Value = -29759.5
$ws.Range("K65").Value = 17896.875
$ws.Range("H65").Value = 11890.083
$ws.Range("J65").Value = 28511.5
$ws.Range("L65").Value = 142557.5
$ws.Range("M65").Value = -14776.875
$ws.Range("I65").Value = 3579.375
$ws.Range("N65").Value = -148797.5
$ws.Range("L68").Value = 39998.5
$ws.Range("J68").Value = 39998.5
$ws.Range("N68").Value = -41496.5
$ws.Range("H68").Value = 38180.453
$ws.Range("H71").Value = 38180.453
$ws.Range("J71").Value = 39998.5
$ws.Range("N71").Value = -127483.5
$ws.Range("L71").Value = 119995.5
$ws.Range("J74").Value = 38983.332
$ws.Range("H74").Value = 37700
$ws.Range("L74").Value = 38983.332
$ws.Range("N74").Value = -40731.332
$ws.Range("N77").Value = -125685.996
$ws.Range("H77").Value = 37700
$ws.Range("L77").Value = 116949.996
$ws.Range("J77").Value = 38983.332
$ws.Range("J81").Value = 199999
$ws.Range("N81").Value = -201995
$ws.Range("H81").Value = 199999
$ws.Range("L81").Value = 199999
$ws.Range("L84").Value = 599997
$ws.Range("J84").Value = 199999
$ws.Range("N84").Value = -609981
$ws.Range("H84").Value = 199999
$ws.Range("M93").Value = -42576.668
$ws.Range("I93").Value = 44448.668
$ws.Range("H93").Value = 45836.375
$ws.Range("K93").Value = 44448.668
$ws.Range("M113").Value = 1216.1539
$ws.Range("I113").Value = 953.8461
$ws.Range("J113").Value = 499.25
$ws.Range("K113").Value = 953.8461
$ws.Range("H113").Value = 846.8823
$ws.Range("L113").Value = 499.25
$ws.Range("N113").Value = -4839.25
$ws.Range("J124").Value = 45000
$ws.Range("L124").Value = 45000
$ws.Range("N124").Value = -49910
$ws.Range("H124").Value = 45000
$ws.Range("L132").Value = 17242.5
$ws.Range("H132").Value = 2046.56
$ws.Range("N132").Value = -22302.5
$ws.Range("M132").Value = -2644.2173
$ws.Range("K132").Value = 5174.2173
$ws.Range("J132").Value = 5747.5
$ws.Range("I132").Value = 1724.7391
$ws.Range("M134").Value = -8274.988799999999
$ws.Range("H134").Value = 3628.29
$ws.Range("I134").Value = 3603.3296
$ws.Range("L134").Value = 11642.0001
$ws.Range("K134").Value = 10809.9888
$ws.Range("N134").Value = -16712.0001
$ws.Range("J134").Value = 3880.6667
$ws.Range("J141").Value = 52995.332
$ws.Range("N141").Value = -63355.332
$ws.Range("L141").Value = 52995.332
$ws.Range("H141").Value = 52995.332

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 138.08333
$ws.Range("L2").Value = 264
$ws.Range("J2").Value = 44
$ws.Range("N2").Value = -490
$ws.Range("L12").Value = 366
$ws.Range("J12").Value = 122
$ws.Range("H12").Value = 225.875
$ws.Range("K12").Value = 1197
$ws.Range("N12").Value = -712
$ws.Range("M12").Value = -1024
$ws.Range("I12").Value = 399
$ws.Range("K14").Value = 379.5
$ws.Range("M14").Value = -206.5
$ws.Range("H14").Value = 126.5
$ws.Range("I14").Value = 126.5
$ws.Range("I38").Value = 247
$ws.Range("N38").Value = -1703.99998
$ws.Range("J38").Value = 336.66666
$ws.Range("K38").Value = 741
$ws.Range("L38").Value = 1009.99998
$ws.Range("M38").Value = -394
$ws.Range("H38").Value = 294.47058
$ws.Range("H70").Value = 17499.75
$ws.Range("N70").Value = -60630
$ws.Range("L70").Value = 60000
$ws.Range("J70").Value = 20000
$ws.Range("N73").Value = -62184
$ws.Range("H73").Value = 17499.75
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 60000
$ws.Range("H119").Value = 12499.5
$ws.Range("J119").Value = 12499.5
$ws.Range("L119").Value = 37498.5
$ws.Range("N119").Value = -47174.5
$ws.Range("M119").ClearContents()
$ws.Range("K119").Value = 0
$ws.Range("I119").Value = 0

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N80").Value = -3762.6666
$ws.Range("M80").Value = -302
$ws.Range("H80").Value = 1650
$ws.Range("J80").Value = 1766.6666
$ws.Range("L80").Value = 1766.6666
$ws.Range("I80").Value = 1300
$ws.Range("K80").Value = 1300
$ws.Range("J83").Value = 1766.6666
$ws.Range("M83").Value = -1508
$ws.Range("L83").Value = 8833.333000000001
$ws.Range("N83").Value = -18817.333
$ws.Range("K83").Value = 6500
$ws.Range("H83").Value = 1650
$ws.Range("I83").Value = 1300
$ws.Range("H107").Value = 680.625
$ws.Range("L107").Value = 975.3333
$ws.Range("N107").Value = -4815.3333
$ws.Range("J107").Value = 975.3333
$ws.Range("M113").Value = -740.3332999999998
$ws.Range("I113").Value = 2910.3333
$ws.Range("K113").Value = 2910.3333
$ws.Range("H113").Value = 2910.3333
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 1342.75
$ws.Range("M132").Value = -495.6364999999996
$ws.Range("K132").Value = 3025.6365
$ws.Range("I132").Value = 1008.5455

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("N22").Value = -2687.842
$ws.Range("I22").Value = 1816.7142
$ws.Range("J22").Value = 2097.842
$ws.Range("M22").Value = -1521.7142
$ws.Range("H22").Value = 2022.1538
$ws.Range("L22").Value = 2097.842
$ws.Range("K22").Value = 1816.7142
$ws.Range("N27").Value = -2311.842
$ws.Range("M27").Value = -1709.7142
$ws.Range("L27").Value = 2097.842
$ws.Range("H27").Value = 2022.1538
$ws.Range("K27").Value = 1816.7142
$ws.Range("J27").Value = 2097.842
$ws.Range("I27").Value = 1816.7142
$ws.Range("K40").Value = 2791.8
$ws.Range("H40").Value = 3076
$ws.Range("I40").Value = 2791.8
$ws.Range("M40").Value = -2655.8
$ws.Range("J55").Value = 922.44446
$ws.Range("H55").Value = 642.94116
$ws.Range("L55").Value = 922.44446
$ws.Range("N55").Value = -1268.44446
$ws.Range("L68").Value = 2749.875
$ws.Range("J68").Value = 2749.875
$ws.Range("N68").Value = -4247.875
$ws.Range("H68").Value = 2532.5264
$ws.Range("H69").Value = 46824.668
$ws.Range("H71").Value = 2532.5264
$ws.Range("J71").Value = 2749.875
$ws.Range("N71").Value = -21237.375
$ws.Range("L71").Value = 13749.375
$ws.Range("H72").Value = 46824.668
$ws.Range("L87").Value = 76999.5
$ws.Range("H87").Value = 76999.5
$ws.Range("J87").Value = 76999.5
$ws.Range("N87").Value = -79245.5
$ws.Range("L88").Value = 56666
$ws.Range("J88").Value = 56666
$ws.Range("N88").Value = -57522
$ws.Range("H88").Value = 56666
$ws.Range("J90").Value = 76999.5
$ws.Range("L90").Value = 230998.5
$ws.Range("N90").Value = -242230.5
$ws.Range("H90").Value = 76999.5
$ws.Range("H91").Value = 56666
$ws.Range("N91").Value = -59630
$ws.Range("J91").Value = 56666
$ws.Range("L91").Value = 56666
$ws.Range("N100").ClearContents()
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1860
$ws.Range("I100").Value = 2401
$ws.Range("K100").Value = 2401
$ws.Range("H100").Value = 2401
$ws.Range("L115").Value = 79000
$ws.Range("H115").Value = 79000
$ws.Range("J115").Value = 79000
$ws.Range("N115").Value = -81350
$ws.Range("H122").Value = 3667.8462
$ws.Range("I122").Value = 2513.2666
$ws.Range("K122").Value = 7539.7998
$ws.Range("M122").Value = -5089.7998
$ws.Range("N136").Value = -14575.3998
$ws.Range("L136").Value = 9475.399800000001
$ws.Range("M136").Value = -1374.3531
$ws.Range("I136").Value = 1308.1177
$ws.Range("H136").Value = 2175.4688
$ws.Range("K136").Value = 3924.3531
$ws.Range("J136").Value = 3158.4666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 657940.2
$ws.Range("L14").Value = 6363
$ws.Range("N14").Value = -6699
$ws.Range("J14").Value = 6363
$ws.Range("N31").Value = -35362.668
$ws.Range("H31").Value = 34666.668
$ws.Range("L31").Value = 34666.668
$ws.Range("J31").Value = 34666.668
$ws.Range("J81").Value = 1300
$ws.Range("N81").Value = -4722
$ws.Range("M81").Value = -13109.5
$ws.Range("I81").Value = 7085.25
$ws.Range("H81").Value = 5507.4546
$ws.Range("K81").Value = 14170.5
$ws.Range("L81").Value = 2600
$ws.Range("L84").Value = 13000
$ws.Range("M84").Value = -65548.5
$ws.Range("J84").Value = 1300
$ws.Range("N84").Value = -23608
$ws.Range("I84").Value = 7085.25
$ws.Range("H84").Value = 5507.4546
$ws.Range("K84").Value = 70852.5
$ws.Range("N92").Value = -59992
$ws.Range("L92").Value = 55000
$ws.Range("J92").Value = 55000
$ws.Range("H92").Value = 55000
$ws.Range("J113").Value = 1127.7858
$ws.Range("H113").Value = 910.2
$ws.Range("L113").Value = 3383.3574
$ws.Range("N113").Value = -7723.357400000001
$ws.Range("H132").Value = 6681.857
$ws.Range("M132").Value = -17384.462
$ws.Range("K132").Value = 19914.462
$ws.Range("I132").Value = 6638.154
$ws.Range("N136").Value = -18462.3339
$ws.Range("L136").Value = 13362.3339
$ws.Range("M136").Value = -39499.66800000001
$ws.Range("I136").Value = 14016.556
$ws.Range("H136").Value = 11625.944
$ws.Range("K136").Value = 42049.66800000001
$ws.Range("J136").Value = 4454.1113
